$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = "18_hazards_to_humans_and_domestic_animals"
$ws.Range("F7").Value = "ppe"
$ws.Range("F8").Value = "ppe"
$ws.Range("F12").Value = "application instructions || env warning - species || env warning - water"
$ws.Range("F21").Value = "pollinator"
$ws.Range("F22").Value = "pollinator"
$ws.Range("F26").Value = "pollinator"
$ws.Range("F29").Value = "pollinator"
$ws.Range("F30").Value = "pollinator"
$ws.Range("F33").Value = "application instructions"
$ws.Range("F34").Value = "134_non-agriculture_use_requirements"
$ws.Range("F35").Value = "application instructions"
$ws.Range("F38").Value = "application instructions"
$ws.Range("F43").Value = "application instructions"
$ws.Range("F46").Value = "application instructions"
$ws.Range("F47").Value = "mixing"
$ws.Range("F94").Value = "154_pesticide_storage"
